# Daily attendance processing - 2025-12-30 06:43:52
# Reorders the "Recorded By" (column G) list so that dnasr281@gmail.com
# is listed first when it appears alongside another recorder name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Text

    if ($val -ne $null -and $val.Contains(", dnasr281@gmail.com")) {
        $parts = $val.Split(",")
        if ($parts.Count -eq 2) {
            $first = $parts[0].Trim()
            $second = $parts[1].Trim()
            if ($second -eq "dnasr281@gmail.com") {
                $cell.Value = $second + ", " + $first
            }
        }
    }
}
